# Apply cryptos list price/volume update (GitHub Actions scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.755.21'
$ws.Range('E2').Value = '  +2.62%  '
$ws.Range('D3').Value = '1.864.23'
$ws.Range('E3').Value = '  +2.40%  '
$ws.Range('D4').Value = "'1.041"
$ws.Range('E4').Value = '  +3.15%  '
$ws.Range('D5').Value = "'324.51"
$ws.Range('E5').Value = '  +3.22%  '
$ws.Range('D6').Value = "'1.037"
$ws.Range('E6').Value = '  +2.78%  '
$ws.Range('D7').Value = "'0.4424"
$ws.Range('E7').Value = '  +2.74%  '
$ws.Range('E8').Value = '  +2.97%  '
$ws.Range('D9').Value = "'0.07465"
$ws.Range('E9').Value = '  +2.68%  '
$ws.Range('D10').Value = "'0.8843"
$ws.Range('E10').Value = '  +1.90%  '
$ws.Range('D11').Value = "'21.65"
$ws.Range('E11').Value = '  +1.61%  '
$ws.Range('D12').Value = '1.877.87'
$ws.Range('E12').Value = '  -13.65%  '
$ws.Range('D13').Value = "'5.555"
$ws.Range('E13').Value = '  +2.34%  '
$ws.Range('D14').Value = "'6.758"
$ws.Range('E14').Value = '  +1.98%  '
$ws.Range('D15').Value = "'0.07237"
$ws.Range('E15').Value = '  +3.85%  '
$ws.Range('D16').Value = "'83.76"
$ws.Range('E16').Value = '  +2.79%  '
$ws.Range('D17').Value = "'1.041"
$ws.Range('E17').Value = '  +2.60%  '
$ws.Range('D18').Value = "'0.000009136"
$ws.Range('E18').Value = '  +2.23%  '
$ws.Range('E19').Value = '  +2.79%  '
$ws.Range('D20').Value = "'15.57"
$ws.Range('E20').Value = '  +1.77%  '
$ws.Range('D21').Value = '27.752.20'
$ws.Range('E21').Value = '  +2.48%  '
$ws.Range('D22').Value = "'5.317"
$ws.Range('E22').Value = '  +2.65%  '
$ws.Range('D23').Value = "'11.31"
$ws.Range('E23').Value = '  +2.78%  '
$ws.Range('D24').Value = "'2.007"
$ws.Range('E24').Value = '  +6.61%  '
$ws.Range('D25').Value = "'158.84"
$ws.Range('E25').Value = '  +2.82%  '
$ws.Range('D26').Value = "'18.85"
$ws.Range('E26').Value = '  +2.77%  '
$ws.Range('D27').Value = "'5.318"
$ws.Range('E27').Value = '  +1.27%  '
$ws.Range('E28').Value = '  +3.85%  '
$ws.Range('D29').Value = "'117.90"
$ws.Range('E29').Value = '  +2.72%  '
$ws.Range('D30').Value = "'0.09058"
$ws.Range('E30').Value = '  +0.97%  '
$ws.Range('D31').Value = "'0.7753"
$ws.Range('E31').Value = '  +2.70%  '
$ws.Range('D32').Value = "'1.214"
$ws.Range('E32').Value = '  +1.74%  '
$ws.Range('D33').Value = "'3.040"
$ws.Range('E33').Value = '  +7.95%  '
$ws.Range('D34').Value = "'4.563"
$ws.Range('E34').Value = '  +2.87%  '
$ws.Range('E35').Value = '  +2.89%  '
$ws.Range('E36').Value = '  +1.33%  '
$ws.Range('D37').Value = "'0.01991"
$ws.Range('E37').Value = '  +3.02%  '
$ws.Range('D38').Value = "'0.05337"
$ws.Range('E38').Value = '  +2.05%  '
$ws.Range('D39').Value = "'2.871"
$ws.Range('E39').Value = '  +4.20%  '
$ws.Range('D40').Value = "'0.5193"
$ws.Range('E40').Value = '  +1.36%  '
$ws.Range('D41').Value = "'0.1692"
$ws.Range('E41').Value = '  +2.23%  '
$ws.Range('D42').Value = "'6.867"
$ws.Range('E42').Value = '  +5.57%  '
$ws.Range('D43').Value = "'8.653"
$ws.Range('E43').Value = '  +3.57%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').Value = "'110.09"
$ws.Range('E44').Value = '  +2.72%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = "'10.61"
$ws.Range('E45').Value = '  +1.92%  '
$ws.Range('D46').Value = "'1.718"
$ws.Range('E46').Value = '  +3.98%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = "'0.06487"
$ws.Range('E47').Value = '  +4.34%  '
$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D48').Value = "'0.4705"
$ws.Range('E48').Value = '  +2.48%  '
$ws.Range('D49').Value = "'1.905"
$ws.Range('E49').Value = '  +2.86%  '
$ws.Range('D50').Value = "'39.81"
$ws.Range('E50').Value = '  +1.19%  '
$ws.Range('D51').Value = "'64.55"
$ws.Range('E51').Value = '  +1.54%  '

# Clear the auto-applied "text" style overrides so cells keep their original (default) styling
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
